# Edit George Garton.xlsx to add matchNo column and two extra match rows
# (Eliminator vs KKR, 52nd vs SRH) ahead of the existing 48th-match row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from the default "Sheet1" to the player's name.
$ws.Name = "George Garton"

# Insert a brand-new column A (matchNo) - shifts teamName..result from A:L to B:M.
$ws.Columns("A:A").Insert()

# Insert two new blank rows above the existing data row (currently row 2),
# pushing it down to row 4.
$ws.Rows("2:3").Insert()

# Make sure numeric-looking values that must be kept as text keep their
# original text representation (they were authored as t="str" cells).
$ws.Columns("E:I").NumberFormat = "@"

# The "states" column is blank for the two new rows, but the cell itself
# should still exist (as an empty text cell) rather than be entirely absent.
$ws.Range("D2:D3").NumberFormat = "@"

# ---- Header row ----
$ws.Range("A1").Value = "matchNo"
$ws.Range("B1").Value = "teamName"
$ws.Range("C1").Value = "batterName"
$ws.Range("D1").Value = "states"
$ws.Range("E1").Value = "runs"
$ws.Range("F1").Value = "balls"
$ws.Range("G1").Value = "fours"
$ws.Range("H1").Value = "sixes"
$ws.Range("I1").Value = "sr"
$ws.Range("J1").Value = "opponentTeamName"
$ws.Range("K1").Value = "venue"
$ws.Range("L1").Value = "date"
$ws.Range("M1").Value = "result"

# ---- Row 2: Eliminator vs Kolkata Knight Riders ----
$ws.Range("A2").Value = "Eliminator"
$ws.Range("B2").Value = "Royal Challengers Bangalore"
$ws.Range("C2").Value = "George Garton"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "0"
$ws.Range("F2").Value = "0"
$ws.Range("G2").Value = "0"
$ws.Range("H2").Value = "0"
$ws.Range("I2").Value = "-"
$ws.Range("J2").Value = "Kolkata Knight Riders"
$ws.Range("K2").Value = "Sharjah"
$ws.Range("L2").Value = "October 11"
$ws.Range("M2").Value = "KKR won by 4 wickets (with 2 balls remaining)"

# ---- Row 3: 52nd match vs Sunrisers Hyderabad ----
$ws.Range("A3").Value = "52nd"
$ws.Range("B3").Value = "Royal Challengers Bangalore"
$ws.Range("C3").Value = "George Garton"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "2"
$ws.Range("F3").Value = "3"
$ws.Range("G3").Value = "0"
$ws.Range("H3").Value = "0"
$ws.Range("I3").Value = "66.66"
$ws.Range("J3").Value = "Sunrisers Hyderabad"
$ws.Range("K3").Value = "Abu Dhabi"
$ws.Range("L3").Value = "October 06"
$ws.Range("M3").Value = "Sunrisers won by 4 runs"

# ---- Row 4: 48th match vs Punjab Kings (the original, pre-existing row) ----
$ws.Range("A4").Value = "48th"
$ws.Range("B4").Value = "Royal Challengers Bangalore"
$ws.Range("C4").Value = "George Garton"
$ws.Range("D4").Value = "b Mohammed Shami"
$ws.Range("E4").Value = "0"
$ws.Range("F4").Value = "1"
$ws.Range("G4").Value = "0"
$ws.Range("H4").Value = "0"
$ws.Range("I4").Value = "0.00"
$ws.Range("J4").Value = "Punjab Kings"
$ws.Range("K4").Value = "Sharjah"
$ws.Range("L4").Value = "October 03"
$ws.Range("M4").Value = "RCB won by 6 runs"
